$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 263.3
$ws.Range("I15").Value = 263.3
$ws.Range("K15").Value = 789.9000000000001
$ws.Range("M15").Value = -620.9000000000001
$ws.Range("H98").Value = 2473.389
$ws.Range("I98").Value = 2560.0881
$ws.Range("J98").Value = 999.5
$ws.Range("K98").Value = 2560.0881
$ws.Range("L98").Value = 999.5
$ws.Range("M98").Value = -1062.0881
$ws.Range("N98").Value = -3995.5
$ws.Range("H111").Value = 1480.1428
$ws.Range("I111").Value = 1265.8
$ws.Range("K111").Value = 3797.4
$ws.Range("M111").Value = -730.3999999999996
$ws.Range("H116").Value = 1447.5
$ws.Range("I116").Value = 1447.5
$ws.Range("K116").Value = 1447.5
$ws.Range("M116").Value = 1994.5
$ws.Range("H122").Value = 2473.389
$ws.Range("I122").Value = 2560.0881
$ws.Range("J122").Value = 999.5
$ws.Range("K122").Value = 7680.2643
$ws.Range("L122").Value = 2998.5
$ws.Range("M122").Value = -5230.2643
$ws.Range("N122").Value = -7898.5
$ws.Range("H125").Value = 2680
$ws.Range("I125").Value = 1133.3334
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 10200.0006
$ws.Range("L125").Value = 45000
$ws.Range("M125").Value = -7740.000599999999
$ws.Range("N125").Value = -49920
$ws.Range("H129").Value = 961
$ws.Range("J129").Value = 961
$ws.Range("L129").Value = 2883
$ws.Range("N129").Value = -12883
$ws.Range("H137").Value = 18519536
$ws.Range("I137").Value = 1007.97437
$ws.Range("K137").Value = 3023.92311
$ws.Range("M137").Value = -473.9231100000002
$ws.Range("H138").Value = 1565.24
$ws.Range("I138").Value = 738.1967
$ws.Range("J138").Value = 2858.8206
$ws.Range("K138").Value = 2214.5901
$ws.Range("L138").Value = 8576.461800000001
$ws.Range("M138").Value = 2925.4099
$ws.Range("N138").Value = -18856.4618
$ws.Range("H141").Value = 1178.3658
$ws.Range("I141").Value = 614.37036
$ws.Range("J141").Value = 2266.0715
$ws.Range("K141").Value = 1843.11108
$ws.Range("L141").Value = 6798.2145
$ws.Range("M141").Value = 3336.88892
$ws.Range("N141").Value = -17158.2145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1644.6666
$ws.Range("I61").Value = 1594.1
$ws.Range("J61").Value = 1897.5
$ws.Range("K61").Value = 1594.1
$ws.Range("L61").Value = 1897.5
$ws.Range("M61").Value = -1382.1
$ws.Range("N61").Value = -2321.5
$ws.Range("H74").Value = 4015.2222
$ws.Range("I74").Value = 797.1111
$ws.Range("J74").Value = 16887.666
$ws.Range("K74").Value = 797.1111
$ws.Range("L74").Value = 16887.666
$ws.Range("M74").Value = 76.88890000000004
$ws.Range("N74").Value = -18635.666
$ws.Range("H77").Value = 4015.2222
$ws.Range("I77").Value = 797.1111
$ws.Range("J77").Value = 16887.666
$ws.Range("K77").Value = 3985.5555
$ws.Range("L77").Value = 84438.33
$ws.Range("M77").Value = 382.4445000000001
$ws.Range("N77").Value = -93174.33
$ws.Range("H132").Value = 28112.895
$ws.Range("I132").Value = 1658.069
$ws.Range("J132").Value = 113356.22
$ws.Range("K132").Value = 4974.207
$ws.Range("L132").Value = 340068.66
$ws.Range("M132").Value = -2444.207
$ws.Range("N132").Value = -345128.66
$ws.Range("H136").Value = 1644.6666
$ws.Range("I136").Value = 1594.1
$ws.Range("J136").Value = 1897.5
$ws.Range("K136").Value = 4782.299999999999
$ws.Range("L136").Value = 5692.5
$ws.Range("M136").Value = -2232.299999999999
$ws.Range("N136").Value = -10792.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 559473.7
$ws.Range("I31").Value = 1346.0769
$ws.Range("J31").Value = 2978026.8
$ws.Range("K31").Value = 1346.0769
$ws.Range("L31").Value = 2978026.8
$ws.Range("M31").Value = -1051.0769
$ws.Range("N31").Value = -2978616.8
$ws.Range("H34").Value = 559473.7
$ws.Range("I34").Value = 1346.0769
$ws.Range("J34").Value = 2978026.8
$ws.Range("K34").Value = 1346.0769
$ws.Range("L34").Value = 2978026.8
$ws.Range("M34").Value = -1144.0769
$ws.Range("N34").Value = -2978430.8
$ws.Range("H58").Value = 1600.6666
$ws.Range("I58").Value = 1151.5
$ws.Range("J58").Value = 2948.1667
$ws.Range("K58").Value = 1151.5
$ws.Range("L58").Value = 2948.1667
$ws.Range("M58").Value = -948.5
$ws.Range("N58").Value = -3354.1667
$ws.Range("H105").Value = 1192.3
$ws.Range("I105").Value = 966.3333
$ws.Range("K105").Value = 966.3333
$ws.Range("M105").Value = 780.6667
$ws.Range("H132").Value = 1175
$ws.Range("I132").Value = 863.1622
$ws.Range("J132").Value = 3482.6
$ws.Range("K132").Value = 2589.4866
$ws.Range("L132").Value = 10447.8
$ws.Range("M132").Value = -59.48660000000018
$ws.Range("N132").Value = -15507.8
$ws.Range("H134").Value = 1254
$ws.Range("I134").Value = 1316.375
$ws.Range("J134").Value = 921.3333
$ws.Range("K134").Value = 3949.125
$ws.Range("L134").Value = 2763.9999
$ws.Range("M134").Value = -1414.125
$ws.Range("N134").Value = -7833.9999
$ws.Range("H136").Value = 1600.6666
$ws.Range("I136").Value = 1151.5
$ws.Range("J136").Value = 2948.1667
$ws.Range("K136").Value = 3454.5
$ws.Range("L136").Value = 8844.500100000001
$ws.Range("M136").Value = -904.5
$ws.Range("N136").Value = -13944.5001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 50000670
$ws.Range("I34").Value = 226.42857
$ws.Range("J34").Value = 76923990
$ws.Range("K34").Value = 679.28571
$ws.Range("L34").Value = 230771970
$ws.Range("M34").Value = -595.28571
$ws.Range("N34").Value = -230772138
$ws.Range("H39").Value = 4409.091
$ws.Range("J39").Value = 4409.091
$ws.Range("L39").Value = 13227.273
$ws.Range("N39").Value = -13815.273
$ws.Range("H55").Value = 476.47058
$ws.Range("J55").Value = 681.8182
$ws.Range("L55").Value = 2045.4546
$ws.Range("N55").Value = -2399.4546
$ws.Range("H76").Value = 2748.2
$ws.Range("I76").Value = 2042
$ws.Range("J76").Value = 3807.5
$ws.Range("K76").Value = 6126
$ws.Range("L76").Value = 11422.5
$ws.Range("M76").Value = -5743
$ws.Range("N76").Value = -12188.5
$ws.Range("H79").Value = 2748.2
$ws.Range("I79").Value = 2042
$ws.Range("J79").Value = 3807.5
$ws.Range("K79").Value = 6126
$ws.Range("L79").Value = 11422.5
$ws.Range("M79").Value = -4800
$ws.Range("N79").Value = -14074.5
$ws.Range("H110").Value = 1900
$ws.Range("I110").Value = 1900
$ws.Range("K110").Value = 5700
$ws.Range("M110").Value = -1610

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3972
$ws.Range("I70").Value = 3802.182
$ws.Range("J70").Value = 4283.3335
$ws.Range("K70").Value = 3802.182
$ws.Range("L70").Value = 4283.3335
$ws.Range("M70").Value = -3532.182
$ws.Range("N70").Value = -4823.3335
$ws.Range("H73").Value = 3972
$ws.Range("I73").Value = 3802.182
$ws.Range("J73").Value = 4283.3335
$ws.Range("K73").Value = 3802.182
$ws.Range("L73").Value = 4283.3335
$ws.Range("M73").Value = -2866.182
$ws.Range("N73").Value = -6155.3335
$ws.Range("H107").Value = 486.72415
$ws.Range("I107").Value = 237.66667
$ws.Range("J107").Value = 894.2727
$ws.Range("K107").Value = 237.66667
$ws.Range("L107").Value = 894.2727
$ws.Range("M107").Value = 1682.33333
$ws.Range("N107").Value = -4734.2727
$ws.Range("H126").Value = 2444.4443
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 1950.4651
$ws.Range("I132").Value = 1894.5385
$ws.Range("J132").Value = 2036
$ws.Range("K132").Value = 5683.6155
$ws.Range("L132").Value = 6108
$ws.Range("M132").Value = -3153.6155
$ws.Range("N132").Value = -11168

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1100.0555
$ws.Range("I82").Value = 898.6087
$ws.Range("J82").Value = 1456.4615
$ws.Range("K82").Value = 898.6087
$ws.Range("L82").Value = 1456.4615
$ws.Range("M82").Value = -537.6087
$ws.Range("N82").Value = -2178.4615
$ws.Range("H85").Value = 1100.0555
$ws.Range("I85").Value = 898.6087
$ws.Range("J85").Value = 1456.4615
$ws.Range("K85").Value = 898.6087
$ws.Range("L85").Value = 1456.4615
$ws.Range("M85").Value = 349.3913
$ws.Range("N85").Value = -3952.4615
$ws.Range("H132").Value = 5382.7417
$ws.Range("I132").Value = 6632.4546
$ws.Range("J132").Value = 2327.889
$ws.Range("K132").Value = 19897.3638
$ws.Range("L132").Value = 6983.667
$ws.Range("M132").Value = -17367.3638
$ws.Range("N132").Value = -12043.667
$ws.Range("H138").Value = 88607
$ws.Range("J138").Value = 88607
$ws.Range("L138").Value = 88607
$ws.Range("N138").Value = -98887

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 30258.5
$ws.Range("I68").Value = 32246
$ws.Range("J68").Value = 28271
$ws.Range("K68").Value = 32246
$ws.Range("L68").Value = 28271
$ws.Range("M68").Value = -31435
$ws.Range("N68").Value = -29893
$ws.Range("H71").Value = 30258.5
$ws.Range("I71").Value = 32246
$ws.Range("J71").Value = 28271
$ws.Range("K71").Value = 96738
$ws.Range("L71").Value = 84813
$ws.Range("M71").Value = -92682
$ws.Range("N71").Value = -92925
$ws.Range("H132").Value = 3497.0173
$ws.Range("I132").Value = 4011.6223
$ws.Range("J132").Value = 1715.6923
$ws.Range("K132").Value = 12034.8669
$ws.Range("L132").Value = 5147.0769
$ws.Range("M132").Value = -9504.866900000001
$ws.Range("N132").Value = -10207.0769
$ws.Range("H136").Value = 4702.537
$ws.Range("I136").Value = 7125.0938
$ws.Range("J136").Value = 1178.8182
$ws.Range("K136").Value = 21375.2814
$ws.Range("L136").Value = 3536.4546
$ws.Range("M136").Value = -18825.2814
$ws.Range("N136").Value = -8636.454600000001
